# Updates the cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are textual prices (may use "." as both thousands and
# decimal separators, e.g. "56.814.29"). Force the cell to Text format before
# assigning so Excel does not reinterpret/round the string as a number, then
# reset the style afterwards so no stray number-format style is left behind.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '56.814.29'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +4.29%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.249.80'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '395.70'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '109.16'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.579'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +5.36%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.246.22'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.05%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '39.23'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0970'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +10.26%  '
$ws.Range('E13').Value = '  +1.72%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.760.01'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.31'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.15'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.255.01'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('E18').Value = '  -3.70%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.66'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '56.716.03'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.27%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.36'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('E22').Value = '  +8.92%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.90'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '294.75'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +7.26%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.11'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.97%  '
$ws.Range('E26').Value = '  -3.32%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '28.20'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('E29').Value = '  -5.27%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.23'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  +2.20%  '
$ws.Range('E34').Value = '  -4.25%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '39.84'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +7.29%  '
$ws.Range('E36').Value = '  -4.38%  '
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '51.51'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.48'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.69%  '
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '136.18'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +4.36%  '
$ws.Range('E43').Value = '  +3.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.89'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.91%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '17.00'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.94'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -5.49%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.278'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.51%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.24'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.14'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.153.03'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('E51').Value = '  -6.13%  '
